$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $oldText"
        return
    }
    # Assign the literal replacement text directly onto the matched range
    # (rather than passing it through Find.Execute's replace argument) so
    # that Word's smart-quote autocorrect doesn't mangle straight
    # apostrophes in the new text.
    $rng.Text = $newText
}

# 1. Update phone number.
Replace-Text "+38 098 514 04 59" "+38 096 707 48 15"

# Remove the now-empty paragraph that used to follow the phone number line.
$phoneParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Phone number:*") {
        $phoneParaIndex = $i
        break
    }
}
if ($phoneParaIndex -gt 0) {
    $nextPara = $d.Paragraphs.Item($phoneParaIndex + 1)
    if ($nextPara.Range.Text.Trim() -eq "") {
        $nextPara.Range.Delete()
    }
}

# 2. English level text tweak: drop the trailing " :(" from "Intermediate :(".
# The surrounding " Pre-" run shares identical formatting with "Intermediate :(",
# so a naive text replace lets the engine coalesce them into a single run.
# Touching Font.Name on the doomed " :(" sub-range (even to its own value)
# before deleting it keeps "Intermediate" as its own distinct run, matching
# the diff (which only rewrites the <w:t> content, not the run split).
$rng = $d.Content
$found = $rng.Find.Execute("Intermediate :(", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    Write-Output "NOT FOUND: Intermediate :("
} else {
    $suffixRng = $d.Range($rng.Start + 12, $rng.End)
    if ($suffixRng.Text -ne " :(") {
        Write-Output "UNEXPECTED SUFFIX: [$($suffixRng.Text)]"
    }
    $suffixRng.Font.Name = $suffixRng.Font.Name
    $suffixRng.Delete()
}

# 3. "(Currently studying)" parenthetical note.
Replace-Text "(Currently studying)" "(I'm studying at the current time)"

# 4. Skills bullet about MIPT coursework (note the leading Cyrillic C in the source text).
Replace-Text "Сurrently studying algorithms and data structures for the course of lectures at MIPT" "I'm studying at the current time algorithms and data structures for the course of lectures at MIPT"

# 5-6-7. Jooble experience bullets.
Replace-Text "parsing data with requests+beautifulsoup" "I parsed data using requests+beautifulsoup libs"
Replace-Text "writing SQL queries to monitor the work of remote employees of the company" "I wrote SQL queries to monitor work of remote employees of the company"
Replace-Text "creating ideas for improving the work of remote employees" "I created a solution for improve the work of remote employees."

# 8-9. Apex Quest Project bullets.
Replace-Text "creating an engine for non-linear street quests with Python/Django and JS/Jquery" "I creating an engine for non-linear street quests using Python/Django and JS/Jquery"
Replace-Text "creating and administering large database" "I create and administer large database"

# 10-11. Commercial project bullets.
Replace-Text "created interactive map using" "I created interactive map using"
Replace-Text "wrote api with Python/Django" "I wrote api using Python/Django"

# 12. Expected salary.
Replace-Text "Expected salary: 800-1000$" "Expected salary: 700-1000$"

Write-Output "edits applied"
